$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column C entirely (it becomes unused / absent from the sheet data)
$ws.Range("C1:C3").ClearContents()

# New data values for A:B and D:E across rows 1-4
$data = @(
    @(1, 22, 21, 6),
    @(2, 35, 11, 5),
    @(3, 67, 54, 9),
    @(4, 86, 75, 0)
)

for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}

# Apply wrap text formatting to D3
$ws.Range("D3").WrapText = $true

# Update selection to E9
$ws.Range("E9").Select()
